# Tribals.docx - "Saving date tribal fees were paid, only works for one
# date currently"
#
# Collapses the old two-blank-fill-in-the-date layout
#     __/__/2015 & ___/__/2015 - Dates Tribal Fees were paid
# into a single merge-style placeholder, matching the `_xxx_` convention
# used for the other fields in this template (_invoice_num_, _subdivision_,
# _reference_num_, _tribe_, _amount_, _trans_ref_num_, ...):
#     _date_paid_ - Dates Tribal Fees were paid
#
# The stray "_GoBack" bookmark that used to sit in front of _trans_ref_num_
# (an artifact of wherever the cursor last was when the file was saved)
# moves up into the newly-edited paragraph, right after the new
# "_date_paid_ " text.

$d = $word.ActiveDocument

# --- locate the region to replace -----------------------------------------
# From the start of the old "__/__/2015" run through the end of
# "..._trans_ref_num_" two paragraphs later. Using Find to locate the
# boundaries (rather than hard-coded character offsets) keeps this robust
# against the unrelated content earlier in the document.

$startRange = $d.Content
$startRange.Find.Execute("__/__/2015", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
if (-not $startRange.Find.Found) {
    throw "Could not find the start of the Tribal Fees date paragraph"
}
$startPos = $startRange.Start

$endRange = $d.Content
$endRange.Find.Execute("trans_ref_num_", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
if (-not $endRange.Find.Found) {
    throw "Could not find the end of the Tribal Consultation Fees paragraph"
}
$endPos = $endRange.End

$target = $d.Range($startPos, $endPos)

# --- build the replacement paragraphs --------------------------------------
# Inserted as raw WordOpenXML so the run/proofErr/bookmark layout matches
# exactly what Word itself produces for a typed "_date_paid_" placeholder.

$body = '<w:p>' + `
            '<w:r><w:t>_</w:t></w:r>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r><w:t>date_paid</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r><w:t>_</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
            '<w:bookmarkEnd w:id="0"/>' + `
            '<w:r><w:t>- Dates Tribal Fees were paid</w:t></w:r>' + `
        '</w:p>' + `
        '<w:p/>' + `
        '<w:p>' + `
            '<w:r><w:t>Tr</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve">ibal Consultation Fees for </w:t></w:r>' + `
            '<w:r><w:t>_</w:t></w:r>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r><w:t>trans_ref_num</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r><w:t>_</w:t></w:r>' + `
        '</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" ' + `
               'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
               '<pkg:xmlData>' + `
                   '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
                       '<w:body>' + $body + '</w:body>' + `
                   '</w:document>' + `
               '</pkg:xmlData>' + `
           '</pkg:part>' + `
       '</pkg:package>'

$target.InsertXML($xml)
